$wb = $excel.ActiveWorkbook

# The two sheets "展览" and "全部类型" share identical rows 2,3,5,6 data,
# and both need their "想去人数" (F column) figures refreshed.
$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("F2").Value = 6577
    $ws.Range("F3").Value = 40
    $ws.Range("F5").Value = 1027
    $ws.Range("F6").Value = 128
}
